$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data (new test case values)
$ws.Range("A2").Value = "fRHry839"
$ws.Range("B2").Value = 231009373
$ws.Range("C2").Value = "cjmmreo86"
$ws.Range("D2").Value = "kPj8&G!7"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "qkXVVxCk"
$ws.Range("G2").Value = "LVNW"
$ws.Range("H2").Value = "Candidate"

# Update row 3 data (new test case values)
$ws.Range("A3").Value = "KyexF639"
$ws.Range("B3").Value = 231009372
$ws.Range("C3").Value = "tyvdrww63"
$ws.Range("D3").Value = "EaJ$%u67"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "ZsiKulkC"
$ws.Range("G3").Value = "Tisc"
$ws.Range("H3").Value = "Candidate"

# Remove row 4 entirely (reduces the data range from A1:H4 to A1:H3)
$ws.Range("A4:H4").Delete()

# Update the selection to match the new used range
[void]$ws.Range("A1:H3").Select()
